# Update column G ("K") values on Sheet1 per regenerated save_data
# (commit: regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 3
    3  = 3
    4  = 1
    5  = 1
    6  = 3
    7  = 1
    8  = 5
    9  = 3
    10 = 3
    11 = 4
    12 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
